# Update the division-fact answers in the worksheet table.
# The table has 20 rows (5 data rows with answers, each followed by
# 3 blank rows) x 5 columns. Each populated cell holds a single run
# with the "N÷N=Q, R" text that needs to be replaced in place; setting
# Cell.Range.Text preserves the existing run/paragraph formatting
# (TimeNewRoman, sz 30) already present in the cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "74÷8=9, 2"
$t.Cell(1, 2).Range.Text = "67÷8=8, 3"
$t.Cell(1, 3).Range.Text = "22÷9=2, 4"
$t.Cell(1, 4).Range.Text = "86÷9=9, 5"
$t.Cell(1, 5).Range.Text = "74÷2=37, 0"

# Row 5
$t.Cell(5, 1).Range.Text = "54÷2=27, 0"
$t.Cell(5, 2).Range.Text = "83÷9=9, 2"
$t.Cell(5, 3).Range.Text = "17÷9=1, 8"
$t.Cell(5, 4).Range.Text = "98÷2=49, 0"
$t.Cell(5, 5).Range.Text = "50÷6=8, 2"

# Row 9
$t.Cell(9, 1).Range.Text = "12÷9=1, 3"
$t.Cell(9, 2).Range.Text = "10÷9=1, 1"
$t.Cell(9, 3).Range.Text = "36÷9=4, 0"
$t.Cell(9, 4).Range.Text = "29÷2=14, 1"
$t.Cell(9, 5).Range.Text = "28÷3=9, 1"

# Row 13
$t.Cell(13, 1).Range.Text = "63÷2=31, 1"
$t.Cell(13, 2).Range.Text = "58÷6=9, 4"
$t.Cell(13, 3).Range.Text = "29÷8=3, 5"
$t.Cell(13, 4).Range.Text = "92÷2=46, 0"
$t.Cell(13, 5).Range.Text = "93÷3=31, 0"

# Row 17
$t.Cell(17, 1).Range.Text = "45÷9=5, 0"
$t.Cell(17, 2).Range.Text = "34÷7=4, 6"
$t.Cell(17, 3).Range.Text = "15÷8=1, 7"
$t.Cell(17, 4).Range.Text = "38÷4=9, 2"
$t.Cell(17, 5).Range.Text = "13÷4=3, 1"
